# Adds two new metric blocks ("accuracy_balanced_mean" / "accuracy_balanced_std")
# to Sheet1, per reviewer feedback. The existing f1_macro_std / f1_micro_std blocks
# are pushed down to make room, keeping their original row order and values intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 6 blank rows right before the old "f1_macro_std" block (row 14).
#    This pushes the existing f1_macro_std block down to 20:25 and the existing
#    f1_micro_std block down to 26:31 (values/types/styles move with the rows).
$ws.Rows("14:19").Insert()

# 2) Insert 6 more blank rows right after the (now shifted) "f1_micro_std" block
#    (26:31), making room for the new "accuracy_balanced_std" rows at 32:37.
$ws.Rows("32:37").Insert()

# 3) The "n_sample" column (B) stores its values as text (e.g. "100", not 100)
#    everywhere else in the sheet, so force text format before writing into the
#    two freshly-inserted (still empty) blocks.
$ws.Range("B14:B19").NumberFormat = "@"
$ws.Range("B32:B37").NumberFormat = "@"

# 4) Re-use the existing bold / bordered / centered label style (column A)
#    for the new metric-name cells, matching the rest of the table.
$ws.Range("A2").Copy()
$ws.Range("A14:A19").PasteSpecial(-4122)
$ws.Range("A32:A37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5) Fill in the new "accuracy_balanced_mean" rows (14:19) and the new
#    "accuracy_balanced_std" rows (32:37). Rows 20:31 already hold the correct,
#    untouched f1_macro_std / f1_micro_std data thanks to the row-insert shift above.

# -- accuracy_balanced_mean --
$ws.Range("A14").Value = "accuracy_balanced_mean"
$ws.Range("B14").Value = "0"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0.435
$ws.Range("A15").Value = "accuracy_balanced_mean"
$ws.Range("B15").Value = "100"
$ws.Range("C15").Value = 0.482
$ws.Range("D15").Value = 0.478
$ws.Range("E15").Value = 0.482
$ws.Range("F15").Value = 0.474
$ws.Range("G15").Value = 0.501
$ws.Range("H15").Value = 0.581
$ws.Range("A16").Value = "accuracy_balanced_mean"
$ws.Range("B16").Value = "500"
$ws.Range("C16").Value = 0.5580000000000001
$ws.Range("D16").Value = 0.583
$ws.Range("E16").Value = 0.625
$ws.Range("F16").Value = 0.599
$ws.Range("G16").Value = 0.674
$ws.Range("H16").Value = 0.788
$ws.Range("A17").Value = "accuracy_balanced_mean"
$ws.Range("B17").Value = "1000"
$ws.Range("C17").Value = 0.5610000000000001
$ws.Range("D17").Value = 0.593
$ws.Range("E17").Value = 0.634
$ws.Range("F17").Value = 0.665
$ws.Range("G17").Value = 0.72
$ws.Range("H17").Value = 0.8080000000000001
$ws.Range("A18").Value = "accuracy_balanced_mean"
$ws.Range("B18").Value = "2500"
$ws.Range("C18").Value = 0.6
$ws.Range("D18").Value = 0.625
$ws.Range("E18").Value = 0.699
$ws.Range("F18").Value = 0.715
$ws.Range("G18").Value = 0.757
$ws.Range("H18").Value = 0.838
$ws.Range("A19").Value = "accuracy_balanced_mean"
$ws.Range("B19").Value = "3188 (all)"
$ws.Range("C19").Value = 0.62
$ws.Range("D19").Value = 0.632
$ws.Range("E19").Value = 0.699
$ws.Range("F19").Value = 0.734
$ws.Range("G19").Value = 0.769
$ws.Range("H19").Value = 0.775

# -- accuracy_balanced_std --
$ws.Range("A32").Value = "accuracy_balanced_std"
$ws.Range("B32").Value = "0"
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("A33").Value = "accuracy_balanced_std"
$ws.Range("B33").Value = "100"
$ws.Range("C33").Value = 0.018
$ws.Range("D33").Value = 0.015
$ws.Range("E33").Value = 0.008999999999999999
$ws.Range("F33").Value = 0.01
$ws.Range("G33").Value = 0.019
$ws.Range("H33").Value = 0.029
$ws.Range("A34").Value = "accuracy_balanced_std"
$ws.Range("B34").Value = "500"
$ws.Range("C34").Value = 0.013
$ws.Range("D34").Value = 0.008
$ws.Range("E34").Value = 0.02
$ws.Range("F34").Value = 0.006
$ws.Range("G34").Value = 0.022
$ws.Range("H34").Value = 0.01
$ws.Range("A35").Value = "accuracy_balanced_std"
$ws.Range("B35").Value = "1000"
$ws.Range("C35").Value = 0.008999999999999999
$ws.Range("D35").Value = 0.014
$ws.Range("E35").Value = 0.026
$ws.Range("F35").Value = 0.024
$ws.Range("G35").Value = 0.021
$ws.Range("H35").Value = 0.013
$ws.Range("A36").Value = "accuracy_balanced_std"
$ws.Range("B36").Value = "2500"
$ws.Range("C36").Value = 0.016
$ws.Range("D36").Value = 0.005
$ws.Range("E36").Value = 0.007
$ws.Range("F36").Value = 0.007
$ws.Range("G36").Value = 0.005
$ws.Range("H36").Value = 0.013
$ws.Range("A37").Value = "accuracy_balanced_std"
$ws.Range("B37").Value = "3188 (all)"
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0.007
$ws.Range("H37").Value = 0.013
